$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.065133174868965
$ws.Range("D2").Value = 1.06784982194883
$ws.Range("E2").Value = 1.077796922916633
$ws.Range("F2").Value = 1.082880924073111
$ws.Range("I2").Value = 1.054010344733571
$ws.Range("J2").Value = 1.070089852595505
$ws.Range("K2").Value = 1.070556580613251
$ws.Range("L2").Value = 1.080477279543
$ws.Range("M2").Value = 1.085547988780456
$ws.Range("N2").Value = 1.026990767719589
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.066185571088467
$ws.Range("D3").Value = 1.068672123750078
$ws.Range("E3").Value = 1.078793744189694
$ws.Range("F3").Value = 1.083858243701293
$ws.Range("I3").Value = 1.054313355121299
$ws.Range("J3").Value = 1.070797244362746
$ws.Range("K3").Value = 1.071194731063045
$ws.Range("L3").Value = 1.08129142470236
$ws.Range("M3").Value = 1.08634363957621
$ws.Range("N3").Value = 1.02723268512326
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.06686697869496
$ws.Range("D4").Value = 1.069204566835808
$ws.Range("E4").Value = 1.079439493556402
$ws.Range("F4").Value = 1.084491336966886
$ws.Range("I4").Value = 1.054508484358153
$ws.Range("J4").Value = 1.07125478531488
$ws.Range("K4").Value = 1.071607362411491
$ws.Range("L4").Value = 1.081818349936037
$ws.Range("M4").Value = 1.08685855573058
$ws.Range("N4").Value = 1.027389023089995
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.067153546670354
$ws.Range("D5").Value = 1.069428490924381
$ws.Range("E5").Value = 1.079711142780055
$ws.Range("F5").Value = 1.084757656596113
$ws.Range("I5").Value = 1.054590291613494
$ws.Range("J5").Value = 1.071447089757504
$ws.Range("K5").Value = 1.071780761476564
$ws.Range("L5").Value = 1.082039897088945
$ws.Range("M5").Value = 1.087075044181071
$ws.Range("N5").Value = 1.027454699680572
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.067201668795076
$ws.Range("D6").Value = 1.069466093723666
$ws.Range("E6").Value = 1.079756764188376
$ws.Range("F6").Value = 1.084802382614469
$ws.Range("I6").Value = 1.054604014203313
$ws.Range("J6").Value = 1.071479375818949
$ws.Range("K6").Value = 1.071809871752404
$ws.Range("L6").Value = 1.082077097440556
$ws.Range("M6").Value = 1.087111394561902
$ws.Range("N6").Value = 1.027465724251946
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.066870807421216
$ws.Range("D7").Value = 1.069207558587497
$ws.Range("E7").Value = 1.079443122655324
$ws.Range("F7").Value = 1.084494894886821
$ws.Range("I7").Value = 1.054509578355512
$ws.Range("J7").Value = 1.071257355076817
$ws.Range("K7").Value = 1.071609679658629
$ws.Range("L7").Value = 1.081821310151598
$ws.Range("M7").Value = 1.086861448391183
$ws.Range("N7").Value = 1.02738990085237
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.065488746413986
$ws.Range("D8").Value = 1.068127647857779
$ws.Range("E8").Value = 1.078133649757486
$ws.Range("F8").Value = 1.083211067917565
$ws.Range("I8").Value = 1.054112942682184
$ws.Range("J8").Value = 1.070328957709404
$ws.Range("K8").Value = 1.070772307190039
$ws.Range("L8").Value = 1.080752398691272
$ws.Range("M8").Value = 1.085816866346267
$ws.Range("N8").Value = 1.027072565812176
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.063056740638085
$ws.Range("D9").Value = 1.066227494630386
$ws.Range("E9").Value = 1.07583188944447
$ws.Range("F9").Value = 1.080954216321601
$ws.Range("I9").Value = 1.053406845206003
$ws.Range("J9").Value = 1.068691581802528
$ws.Range("K9").Value = 1.069294520013801
$ws.Range("L9").Value = 1.078869780246803
$ws.Range("M9").Value = 1.083976799320353
$ws.Range("N9").Value = 1.026511871991634
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.061437683839561
$ws.Range("D10").Value = 1.064962651530384
$ws.Range("E10").Value = 1.074301260269477
$ws.Range("F10").Value = 1.079453339649692
$ws.Range("I10").Value = 1.052931309068522
$ws.Range("J10").Value = 1.067599074735219
$ws.Range("K10").Value = 1.068307864332126
$ws.Range("L10").Value = 1.077615371899741
$ws.Range("M10").Value = 1.082750545399884
$ws.Range("N10").Value = 1.026137079206739
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.06073715867734
$ws.Range("D11").Value = 1.064415427722463
$ws.Range("E11").Value = 1.073639409955049
$ws.Range("F11").Value = 1.078804328978411
$ws.Range("I11").Value = 1.052724260949797
$ws.Range("J11").Value = 1.06712579451825
$ws.Range("K11").Value = 1.06788029121402
$ws.Range("L11").Value = 1.077072365934104
$ws.Range("M11").Value = 1.082219680499338
$ws.Range("N11").Value = 1.025974556903842
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.060477033320532
$ws.Range("D12").Value = 1.064212234578245
$ws.Range("E12").Value = 1.07339370851773
$ws.Range("F12").Value = 1.078563390301561
$ws.Range("I12").Value = 1.052647183377862
$ws.Range("J12").Value = 1.066949964953118
$ws.Range("K12").Value = 1.067721420239427
$ws.Range("L12").Value = 1.076870694100805
$ws.Range("M12").Value = 1.082022510926874
$ws.Range("N12").Value = 1.025914153910677
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.060532827429607
$ws.Range("D13").Value = 1.064255817010563
$ws.Range("E13").Value = 1.073446405998223
$ws.Range("F13").Value = 1.078615066444835
$ws.Range("I13").Value = 1.052663724498116
$ws.Range("J13").Value = 1.06698768244247
$ws.Range("K13").Value = 1.06775550091519
$ws.Range("L13").Value = 1.07691395226431
$ws.Range("M13").Value = 1.08206480367762
$ws.Range("N13").Value = 1.025927112136114
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.060715654983198
$ws.Range("D14").Value = 1.064398630278285
$ws.Range("E14").Value = 1.073619097353197
$ws.Range("F14").Value = 1.078784410199132
$ws.Range("I14").Value = 1.05271789317397
$ws.Range("J14").Value = 1.067111261046207
$ws.Range("K14").Value = 1.067867159932876
$ws.Range("L14").Value = 1.077055695169155
$ws.Range("M14").Value = 1.082203382044863
$ws.Range("N14").Value = 1.025969564685497
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.060828311795951
$ws.Range("D15").Value = 1.06448663156436
$ws.Range("E15").Value = 1.073725516667526
$ws.Range("F15").Value = 1.078888766092791
$ws.Range("I15").Value = 1.052751245673937
$ws.Range("J15").Value = 1.067187397646653
$ws.Range("K15").Value = 1.067935949958867
$ws.Range("L15").Value = 1.077143030948417
$ws.Range("M15").Value = 1.082288767059198
$ws.Range("N15").Value = 1.025995716475896
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.061484186719395
$ws.Range("D16").Value = 1.064998978765319
$ws.Range("E16").Value = 1.074345204583612
$ws.Range("F16").Value = 1.079496430930331
$ws.Range("I16").Value = 1.052945026219951
$ws.Range("J16").Value = 1.067630480241407
$ws.Range("K16").Value = 1.068336233732548
$ws.Range("L16").Value = 1.077651412859141
$ws.Range("M16").Value = 1.082785779546913
$ws.Range("N16").Value = 1.0261478603576
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.061895743659794
$ws.Range("D17").Value = 1.065320484771423
$ws.Range("E17").Value = 1.074734165828516
$ws.Range("F17").Value = 1.079877838844873
$ws.Range("I17").Value = 1.053066275192517
$ws.Range("J17").Value = 1.0679083563463
$ws.Range("K17").Value = 1.068587229239603
$ws.Range("L17").Value = 1.077970350927459
$ws.Range("M17").Value = 1.083097572520903
$ws.Range("N17").Value = 1.026243233585604
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.062135849671347
$ws.Range("D18").Value = 1.065508058269898
$ws.Range("E18").Value = 1.074961129276596
$ws.Range("F18").Value = 1.080100392372671
$ws.Range("I18").Value = 1.053136887864516
$ws.Range("J18").Value = 1.068070415766163
$ws.Range("K18").Value = 1.068733597380013
$ws.Range("L18").Value = 1.078156397615862
$ws.Range("M18").Value = 1.083279446771023
$ws.Range("N18").Value = 1.026298840527049
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.062217728398974
$ws.Range("D19").Value = 1.065572023506799
$ws.Range("E19").Value = 1.07503853301118
$ws.Range("F19").Value = 1.080176291725681
$ws.Range("I19").Value = 1.053160946328487
$ws.Range("J19").Value = 1.068125670262926
$ws.Range("K19").Value = 1.06878349947109
$ws.Range("L19").Value = 1.078219837342599
$ws.Range("M19").Value = 1.083341462984437
$ws.Range("N19").Value = 1.026317797207359
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.061851582134538
$ws.Range("D20").Value = 1.065285985612651
$ws.Range("E20").Value = 1.074692424764096
$ws.Range("F20").Value = 1.079836908614299
$ws.Range("I20").Value = 1.053053277680683
$ws.Range("J20").Value = 1.067878545039991
$ws.Range("K20").Value = 1.068560303236978
$ws.Range("L20").Value = 1.077936130282901
$ws.Range("M20").Value = 1.083064118976207
$ws.Range("N20").Value = 1.026233003281586
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.060661814591348
$ws.Range("D21").Value = 1.064356573387797
$ws.Range("E21").Value = 1.073568240204195
$ws.Range("F21").Value = 1.078734539012763
$ws.Range("I21").Value = 1.052701946558561
$ws.Range("J21").Value = 1.067074871113319
$ws.Range("K21").Value = 1.067834280546719
$ws.Range("L21").Value = 1.077013954728586
$ws.Range("M21").Value = 1.082162573687679
$ws.Range("N21").Value = 1.025957064431113
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.059914228553964
$ws.Range("D22").Value = 1.063772621254322
$ws.Range("E22").Value = 1.072862226684549
$ws.Range("F22").Value = 1.078042204109751
$ws.Range("I22").Value = 1.052480063358886
$ws.Range("J22").Value = 1.066569383229427
$ws.Range("K22").Value = 1.067377504823172
$ws.Range("L22").Value = 1.076434289884393
$ws.Range("M22").Value = 1.081595836702641
$ws.Range("N22").Value = 1.025783368509364
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.060310493629411
$ws.Range("D23").Value = 1.064082146558836
$ws.Range("E23").Value = 1.073236421186766
$ws.Range("F23").Value = 1.078409150834191
$ws.Range("I23").Value = 1.052597781344053
$ws.Range("J23").Value = 1.066837369418399
$ws.Range("K23").Value = 1.067619678206301
$ws.Range("L23").Value = 1.076741567395833
$ws.Range("M23").Value = 1.081896264994063
$ws.Range("N23").Value = 1.025875467082745
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.061871536666716
$ws.Range("D24").Value = 1.065301574160091
$ws.Range("E24").Value = 1.074711285479601
$ws.Range("F24").Value = 1.079855402961851
$ws.Range("I24").Value = 1.053059151036074
$ws.Range("J24").Value = 1.067892015551005
$ws.Range("K24").Value = 1.068572470041418
$ws.Range("L24").Value = 1.077951593071242
$ws.Range("M24").Value = 1.083079235160058
$ws.Range("N24").Value = 1.026237625985607
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.063685071961592
$ws.Range("D25").Value = 1.066718394147402
$ws.Range("E25").Value = 1.076426269604502
$ws.Range("F25").Value = 1.081537019647832
$ws.Range("I25").Value = 1.053590236584141
$ws.Range("J25").Value = 1.069115047810479
$ws.Range("K25").Value = 1.069676824075884
$ws.Range("L25").Value = 1.079356366872604
$ws.Range("M25").Value = 1.084452423108938
$ws.Range("N25").Value = 1.026657001786422
